# Correction in sa algorithm and 746 logs
# Update the Fitness column (C) values in Sheet1 to reflect the corrected
# simulated-annealing run log values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-16  -> 7685
$ws.Range("C2:C16").Value = 7685

# Rows 17-52 -> 7660
$ws.Range("C17:C52").Value = 7660

# Rows 53-63 -> 7318
$ws.Range("C53:C63").Value = 7318

# Rows 64-68 -> 7312
$ws.Range("C64:C68").Value = 7312

# Rows 69-143 -> 7310
$ws.Range("C69:C143").Value = 7310

# Rows 144-170 -> 7293
$ws.Range("C144:C170").Value = 7293
